$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Package (column D) values for rows 30-62 to reflect the package refactor
# fr.pds.isintheair.crmtab.tlacouque.uc.admin.ref.customer.* -> fr.pds.isintheair.crmtab.*
# Assignments are ordered to match the new-package introduction order of the source commit.
$ws.Range("D45").Value = "fr.pds.isintheair.crmtab.controller.broadcastreceiver"
$ws.Range("D51:D52").Value = "fr.pds.isintheair.crmtab.controller.adapter"
$ws.Range("D46").Value = "fr.pds.isintheair.crmtab.helper"
$ws.Range("D61:D62").Value = "fr.pds.isintheair.crmtab.helper"
$ws.Range("D32").Value = "fr.pds.isintheair.crmtab.model.asynctask"
$ws.Range("D30:D31").Value = "fr.pds.isintheair.crmtab.model.entity"
$ws.Range("D33:D44").Value = "fr.pds.isintheair.crmtab.model.entity"
$ws.Range("D47").Value = "fr.pds.isintheair.crmtab.model.rest.service"
$ws.Range("D48:D49").Value = "fr.pds.isintheair.crmtab.model.rest"
$ws.Range("D50").Value = "fr.pds.isintheair.crmtab.view.activity"
$ws.Range("D53:D60").Value = "fr.pds.isintheair.crmtab.view.fragment"

# Update the active selection/view position to match the new commit state
$ws.Range("G53").Select()
